# "Added CDS All studies testcase": the SamplesTab query (StatQuery column,
# row 3) is simplified to drop the Tumor / Analyte Type columns so it
# returns every sample for the study instead of a per-tumor-status subset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs001524'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

# Update the StatQuery cell for the SamplesTab row (row 3, column B) with the
# simplified query (no Tumor / Analyte Type columns)
$ws.Range("B3").Value = $newQuery

# Leave the selection where the editor left it when saving the workbook
$ws.Range("B3").Select()
